$d = $word.ActiveDocument

$replacements = @(
    @("2024-09-14 Saturday", "2024-09-15 Sunday"),
    @("27÷2=", "89÷2="),
    @("28÷2=", "14÷8="),
    @("82÷3=", "17÷7="),
    @("88÷2=", "91÷8="),
    @("31÷4=", "32÷8="),
    @("36÷9=", "60÷5="),
    @("31÷7=", "46÷5="),
    @("29÷2=", "67÷4="),
    @("72÷7=", "59÷5="),
    @("83÷4=", "14÷3="),
    @("72÷2=", "63÷2="),
    @("63÷8=", "62÷4="),
    @("41÷6=", "61÷4="),
    @("34÷9=", "44÷4="),
    @("46÷7=", "39÷2="),
    @("33÷8=", "97÷5="),
    @("61÷3=", "39÷4="),
    @("41÷8=", "92÷3="),
    @("40÷7=", "80÷5="),
    @("26÷4=", "41÷2="),
    @("63÷4=", "70÷6="),
    @("96÷9=", "89÷2="),
    @("21÷5=", "34÷7="),
    @("66÷3=", "90÷4="),
    @("75÷5=", "93÷3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
